# Update a handful of imputed values in the RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -14.12869999999999
$ws.Range("E4").Value = 13.4468

$ws.Range("E5").Value = 13.17239999999999

$ws.Range("C7").Value = -12.0011

$ws.Range("E8").Value = 14.19989999999999

$ws.Range("C16").Value = -11.8704
$ws.Range("E16").Value = 12.9125
